$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually differ between rows 2, 3 and 4 (others are identical
# across these rows, so a cyclic rotation leaves them visually unchanged):
#   A, B, D, E, F, G, H, Q, R, AI, AN, AO
$cols = @("A","B","D","E","F","G","H","Q","R","AI","AN","AO")

# Capture the current ("before") values for rows 2, 3 and 4 so we can rotate
# them without clobbering data we still need to read.
$row2 = @{}
$row3 = @{}
$row4 = @{}
foreach ($col in $cols) {
    $row2[$col] = $ws.Range("${col}2").Value()
    $row3[$col] = $ws.Range("${col}3").Value()
    $row4[$col] = $ws.Range("${col}4").Value()
}

# Apply the rotation: old row 2 -> row 3, old row 3 -> row 4, old row 4 -> row 2
foreach ($col in $cols) {
    $ws.Range("${col}2").Value = $row4[$col]
    $ws.Range("${col}3").Value = $row2[$col]
    $ws.Range("${col}4").Value = $row3[$col]
}
